# Updates the cryptocurrency price (column D) and 1h volume change (column E)
# values on the active worksheet to match the latest GitHub Actions scrape.
# Rows are keyed by their sheet row number (header is row 1, data starts row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 2; D = '26.344.45'; E = '  +0.99%  ' },
    @{ Row = 3; D = '1.666.44'; E = '  +0.87%  ' },
    @{ Row = 4; D = '1.011'; E = '  +0.95%  ' },
    @{ Row = 5; D = '219.26'; E = '  +0.88%  ' },
    @{ Row = 6; D = '0.5343'; E = '  +1.51%  ' },
    @{ Row = 7; D = '1.011'; E = $null },
    @{ Row = 8; D = $null; E = '  +2.52%  ' },
    @{ Row = 9; D = '0.06396'; E = '  +1.16%  ' },
    @{ Row = 10; D = '20.86'; E = '  +2.44%  ' },
    @{ Row = 11; D = '0.07851'; E = '  +0.82%  ' },
    @{ Row = 12; D = '4.566'; E = '  +1.33%  ' },
    @{ Row = 13; D = '1.662.77'; E = '  +0.58%  ' },
    @{ Row = 14; D = '1.895.72'; E = '  +0.91%  ' },
    @{ Row = 15; D = '0.5537'; E = '  +0.82%  ' },
    @{ Row = 16; D = '0.0₅8196'; E = '  -0.20%  ' },
    @{ Row = 17; D = '65.91'; E = '  +0.59%  ' },
    @{ Row = 18; D = '26.370.32'; E = '  +1.08%  ' },
    @{ Row = 19; D = $null; E = '  +0.88%  ' },
    @{ Row = 20; D = '4.682'; E = '  +2.46%  ' },
    @{ Row = 21; D = '193.91'; E = '  +1.68%  ' },
    @{ Row = 22; D = '10.30'; E = '  +2.28%  ' },
    @{ Row = 23; D = '6.045'; E = '  +0.18%  ' },
    @{ Row = 24; D = '1.012'; E = '  +0.91%  ' },
    @{ Row = 25; D = '146.45'; E = '  +2.09%  ' },
    @{ Row = 26; D = '0.1230'; E = '  -0.54%  ' },
    @{ Row = 27; D = '7.214'; E = '  -0.22%  ' },
    @{ Row = 28; D = '16.12'; E = '  +0.45%  ' },
    @{ Row = 29; D = '1.500'; E = '  +4.71%  ' },
    @{ Row = 30; D = '0.05869'; E = '  +0.65%  ' },
    @{ Row = 31; D = '1.285'; E = '  +0.96%  ' },
    @{ Row = 32; D = '3.619'; E = '  +2.01%  ' },
    @{ Row = 33; D = '3.286'; E = '  +0.75%  ' },
    @{ Row = 34; D = '1.606'; E = '  +1.37%  ' },
    @{ Row = 35; D = '0.9699'; E = '  +2.53%  ' },
    @{ Row = 36; D = '2.829'; E = '  +1.75%  ' },
    @{ Row = 37; D = '2.424'; E = '  +0.59%  ' },
    @{ Row = 38; D = '0.5828'; E = '  +1.72%  ' },
    @{ Row = 39; D = '0.01602'; E = '  -0.45%  ' },
    @{ Row = 40; D = '0.8633'; E = '  +2.50%  ' },
    @{ Row = 41; D = '1.064.75'; E = '  +3.47%  ' },
    @{ Row = 42; D = '5.841'; E = '  +1.71%  ' },
    @{ Row = 43; D = $null; E = '  +0.91%  ' },
    @{ Row = 44; D = '104.65'; E = '  +0.54%  ' },
    @{ Row = 45; D = '1.807.18'; E = '  +0.73%  ' },
    @{ Row = 46; D = '57.88'; E = '  +1.43%  ' },
    @{ Row = 47; D = $null; E = '  -4.82%  ' },
    @{ Row = 48; D = '1.012'; E = $null },
    @{ Row = 49; D = '0.4390'; E = $null },
    @{ Row = 50; D = '8.013'; E = '  +2.58%  ' },
    @{ Row = 51; D = $null; E = '  +0.46%  ' }
)

function Set-TextCellValue($ws, $rowNum, $colLetter, $newValue) {
    if ($null -eq $newValue) {
        return
    }
    $cell = $ws.Range($colLetter + $rowNum)
    # Force a text number format before assignment so Excel does not
    # reinterpret dotted numeric-looking strings (e.g. "1.011") as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    # Restore the default cell style so we don't leave a stray text format
    # applied to a cell that originally had no explicit style.
    $cell.Style = "Normal"
}

foreach ($row in $rowData) {
    Set-TextCellValue $ws $row.Row "D" $row.D
    Set-TextCellValue $ws $row.Row "E" $row.E
}
